$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{ Row=2; D="27.866.52"; E="  -0.31%  " },
    @{ Row=3; D="1.624.30"; E="  -0.97%  " },
    @{ Row=4; E="  -0.04%  " },
    @{ Row=5; D="211.05"; DForceText=$true; E="  -0.94%  " },
    @{ Row=6; E="  -0.24%  " },
    @{ Row=7; D="1.00"; DForceText=$true; E="  -0.05%  " },
    @{ Row=8; D="23.36"; DForceText=$true },
    @{ Row=9; E="  -1.74%  " },
    @{ Row=10; E="  -0.61%  " },
    @{ Row=11; E="  -0.32%  " },
    @{ Row=12; D="1.855.88"; E="  -0.92%  " },
    @{ Row=13; D="1.618.58"; E="  -1.30%  " },
    @{ Row=14; E="  -1.86%  " },
    @{ Row=15; D="0.560"; DForceText=$true; E="  -2.35%  " },
    @{ Row=16; D="65.29"; DForceText=$true },
    @{ Row=17; D="27.854.52"; E="  -0.36%  " },
    @{ Row=18; D="228.99"; DForceText=$true; E="  -2.10%  " },
    @{ Row=19; E="  -0.56%  " },
    @{ Row=20; E="  +0.11%  " },
    @{ Row=21; E="  -0.12%  " },
    @{ Row=22; E="  -1.29%  " },
    @{ Row=23; D="10.08"; DForceText=$true; E="  -6.20%  " },
    @{ Row=24; E="  -2.43%  " },
    @{ Row=25; D="155.06"; DForceText=$true; E="  +2.56%  " },
    @{ Row=26; E="  -1.19%  " },
    @{ Row=27; E="  -0.22%  " },
    @{ Row=28; D="15.48"; DForceText=$true; E="  -1.31%  " },
    @{ Row=29; E="  -0.01%  " },
    @{ Row=30; E="  -0.75%  " },
    @{ Row=31; E="  -0.59%  " },
    @{ Row=32; E="  +1.91%  " },
    @{ Row=33; E="  -1.66%  " },
    @{ Row=34; D="1.392.28"; E="  -1.35%  " },
    @{ Row=35; E="  -0.32%  " },
    @{ Row=36; D="1.00"; DForceText=$true; E="  +10.97%  " },
    @{ Row=37; E="  -1.01%  " },
    @{ Row=38; E="  +0.38%  " },
    @{ Row=39; E="  -0.51%  " },
    @{ Row=40; E="  -3.12%  " },
    @{ Row=41; E="  -0.07%  " },
    @{ Row=42; E="  -0.13%  " },
    @{ Row=44; E="  -3.51%  " },
    @{ Row=45; D="65.61"; DForceText=$true; E="  -1.39%  " },
    @{ Row=46; E="  -0.89%  " },
    @{ Row=47; D="2.16"; DForceText=$true; E="  -2.08%  " },
    @{ Row=48; D="87.86"; DForceText=$true; E="  -0.06%  " },
    @{ Row=49; E="  -1.36%  " },
    @{ Row=50; E="  +0.94%  " },
    @{ Row=51; E="  -0.49%  " }
)

foreach ($chg in $changes) {
    $row = $chg.Row
    if ($chg.ContainsKey('D')) {
        $cell = $ws.Cells.Item($row, 4)
        if ($chg.ContainsKey('DForceText')) {
            $cell.NumberFormat = "@"
            $cell.Value = $chg.D
            $cell.Style = "Normal"
        } else {
            $cell.Value = $chg.D
        }
    }
    if ($chg.ContainsKey('E')) {
        $ws.Cells.Item($row, 5).Value = $chg.E
    }
}
